$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: only the ShipmentTracking (P) value is cleared; ExpectedRate (Q) and
# Result (R) are left exactly as they were.
$ws.Range("P2").ClearContents()

# Rows 3-25: ShipmentTracking (P), ExpectedRate (Q) and Result (R) are all
# cleared. P is left as a fully empty cell (no cached style), while Q and R
# keep an explicit (but default) cell record, so re-apply the Normal style
# after clearing so the cell node is retained with no style attribute.
$ws.Range("P3:P25").ClearContents()
$ws.Range("Q3:R25").ClearContents()
$ws.Range("Q3:R25").Style = "Normal"

# Row 26: ShipmentTracking (P), ExpectedRate (Q) and Result (R) are cleared
# too, but this time the cells keep the bordered formatting already used by
# their row neighbours (O26/S26), with Result (R) additionally centered like
# the rest of the R column. Resetting to the Normal style first drops the
# column-level defaults (e.g. Q's number format, R's centering) so the
# border formatting applied below lines up with the rest of the sheet
# instead of spawning a one-off style.
$ws.Range("P26:R26").ClearContents()
$ws.Range("P26:R26").Style = "Normal"
$ws.Range("P26").Borders.ColorIndex = 1
$ws.Range("Q26").Borders.ColorIndex = 1
$ws.Range("R26").Borders.ColorIndex = 1
$ws.Range("R26").HorizontalAlignment = -4108
